# Time-tracker update: record a "Fees paid" entry on the April and May
# sheets, and leave the April sheet as the active tab (matching the
# author's last-saved view).

$wb = $excel.ActiveWorkbook

# --- April sheet: row 7 gets a "Fees paid" entry for 200 -----------------
$april = $wb.Worksheets.Item("April")
$april.Range("A7").Value = 45795
$april.Range("B7").Value = "Fees paid"
$april.Range("E7").Value = 200

# --- May sheet: row 8 gets a "Fees paid" entry for 300 --------------------
$may = $wb.Worksheets.Item("May")
$may.Range("A8").Value = 45827
$may.Range("B8").Value = "Fees paid"
$may.Range("E8").Value = 300

# --- Selection / active tab bookkeeping -----------------------------------
# (May is activated/selected first, then April last, so April ends up as
# the workbook's active/selected tab - matching the author's last save.)
$may.Activate() | Out-Null
$may.Range("E9").Select() | Out-Null

$april.Activate() | Out-Null
$april.Range("D4").Select() | Out-Null
